# Actualización desde MV -datos-
# Update the "Variación anual" monthly table:
#   - corrects the last existing row (01-05-2021, row 54)
#   - appends a new row for 01-06-2021 (row 55)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct existing row 54 (Serie = 01-05-2021) ---
$ws.Range("B54").Value = 5.6
$ws.Range("C54").Value = 1.9
$ws.Range("D54").Value = 6.1

# --- Append new row 55 (Serie = 01-06-2021) ---
# The date-like label must stay a literal text value (as it is for every
# other "Serie" cell in column A) instead of being auto-converted to a
# date serial number. Writing it through a formula and then converting
# that formula to a static value via Copy/PasteSpecial(values) keeps the
# cell's original "General" style untouched while storing it as text.
$ws.Range("A55").Formula = '="01-06-2021"'
$ws.Range("A55").Copy()
$ws.Range("A55").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("B55").Value = 5.9
$ws.Range("C55").Value = 2
$ws.Range("D55").Value = 6.3
